$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert new column "thr_as" after "thr_gs" (currently column 27), before "low_speed_fix" (28)
$newCol = $tbl.ListColumns.Add(28)
$newCol.Name = "thr_as"

Write-Host $tbl.ListColumns.Count
for ($i = 1; $i -le $tbl.ListColumns.Count; $i++) {
    $col = $tbl.ListColumns.Item($i)
    Write-Host "$i : $($col.Name)"
}
Write-Host $tbl.Range.Address()
